# Update "LP1912" sheet (main schedule for line 141 - LP1912)
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 06:52:53"
$ws1.Range("A3").Value = "Total filas: 51"

$sheet1Rows = @(
    @("05:42:22", "05:43", "14_ABASTO", 1, "LP1912"),
    @("05:42:22", "05:52", "17_ROMERO", 10, "LP1912"),
    @("05:42:22", "06:01", "16_SANTA ANA", 19, "LP1912"),
    @("05:42:22", "06:04", "10_OLMOS", 22, "LP1912"),
    @("05:42:22", "06:11", "215A_EL PATO", 29, "LP1912"),
    @("05:42:22", "06:24", "11_ETCHEVERRY", 42, "LP1912"),
    @("05:42:22", "06:27", "23_HERNANDEZ", 45, "LP1912"),
    @("05:42:22", "06:31", "16_SANTA ANA", 49, "LP1912"),
    @("05:42:22", "06:31", "17X38_ROMERO", 49, "LP1912"),
    @("05:42:22", "06:39", "225_C ROCA-H SUR", 57, "LP1912"),
    @("06:33:46", "06:50", "215A_EL PATO", 17, "LP1912"),
    @("06:16:15", "06:51", "215A_EL PATO", 35, "LP1912"),
    @("06:52:52", "06:52", "215A_EL PATO", 0, "LP1912"),
    @("05:42:22", "06:54", "14_ABASTO", 72, "LP1912"),
    @("06:52:52", "06:55", "14_ABASTO", 3, "LP1912"),
    @("06:45:50", "07:00", "16_SANTA ANA", 15, "LP1912"),
    @("05:57:08", "07:01", "16_SANTA ANA", 64, "LP1912"),
    @("06:45:50", "07:03", "225_GOMEZ", 18, "LP1912"),
    @("05:42:22", "07:04", "225_GOMEZ", 82, "LP1912"),
    @("06:33:46", "07:06", "215C_EL PATO", 33, "LP1912"),
    @("05:42:22", "07:07", "215C_EL PATO", 85, "LP1912"),
    @("06:33:46", "07:13", "14X44_ABASTO", 40, "LP1912"),
    @("05:42:22", "07:14", "14X44_ABASTO", 92, "LP1912"),
    @("06:33:46", "07:15", "16_SANTA ANA", 42, "LP1912"),
    @("06:33:46", "07:20", "215A_EL PATO", 47, "LP1912"),
    @("05:42:22", "07:21", "215A_EL PATO", 99, "LP1912"),
    @("06:52:52", "07:28", "14_ABASTO", 36, "LP1912"),
    @("05:57:08", "07:29", "14_ABASTO", 92, "LP1912"),
    @("05:42:22", "07:33", "23_HERNANDEZ", 111, "LP1912"),
    @("05:57:08", "07:34", "23_HERNANDEZ", 97, "LP1912"),
    @("06:45:50", "07:35", "17X38_ROMERO", 50, "LP1912"),
    @("05:42:22", "07:36", "27_EL RETIRO", 63, "LP1912"),
    @("05:42:22", "07:36", "17X38_ROMERO", 114, "LP1912"),
    @("05:42:22", "07:37", "27_EL RETIRO", 115, "LP1912"),
    @("06:33:46", "07:43", "10_OLMOS", 70, "LP1912"),
    @("05:57:08", "07:44", "10_OLMOS", 107, "LP1912"),
    @("05:57:08", "07:51", "15_ABASTO", 114, "LP1912"),
    @("06:16:15", "07:58", "23_HERNANDEZ", 102, "LP1912"),
    @("06:33:46", "07:59", "11_ETCHEVERRY", 86, "LP1912"),
    @("06:16:15", "08:00", "11_ETCHEVERRY", 104, "LP1912"),
    @("06:33:46", "08:00", "23_HERNANDEZ", 87, "LP1912"),
    @("06:45:50", "08:01", "23_HERNANDEZ", 76, "LP1912"),
    @("06:52:52", "08:02", "23_HERNANDEZ", 70, "LP1912"),
    @("06:16:15", "08:03", "17_ROMERO", 107, "LP1912"),
    @("06:33:46", "08:12", "10_OLMOS", 99, "LP1912"),
    @("06:16:15", "08:15", "17_ROMERO", 119, "LP1912"),
    @("06:33:46", "08:26", "15X38_ABASTO", 113, "LP1912"),
    @("06:33:46", "08:27", "84_COLONIA URQUIZA-ESC 49", 114, "LP1912"),
    @("06:45:50", "08:29", "14_ABASTO", 104, "LP1912"),
    @("06:33:46", "08:31", "16_P MOR-SANTA ANA", 118, "LP1912"),
    @("06:45:50", "08:38", "215C_EL PATO", 113, "LP1912"),
)

$r = 6
foreach ($row in $sheet1Rows) {
    $ws1.Cells.Item($r, 1).Value = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $ws1.Cells.Item($r, 3).Value = $row[2]
    $ws1.Cells.Item($r, 4).Value = $row[3]
    $ws1.Cells.Item($r, 5).Value = $row[4]
    $r++
}

# Update "LP1912-215" sheet (215-filtered subset)
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: 06:52:53"
$ws2.Range("A3").Value = "Total filas: 9"

$sheet2Rows = @(
    @("05:42:22", "06:11", "215A_EL PATO", 29, "LP1912"),
    @("06:33:46", "06:50", "215A_EL PATO", 17, "LP1912"),
    @("06:16:15", "06:51", "215A_EL PATO", 35, "LP1912"),
    @("06:52:52", "06:52", "215A_EL PATO", 0, "LP1912"),
    @("06:33:46", "07:06", "215C_EL PATO", 33, "LP1912"),
    @("05:42:22", "07:07", "215C_EL PATO", 85, "LP1912"),
    @("06:33:46", "07:20", "215A_EL PATO", 47, "LP1912"),
    @("05:42:22", "07:21", "215A_EL PATO", 99, "LP1912"),
    @("06:45:50", "08:38", "215C_EL PATO", 113, "LP1912"),
)

$r = 6
foreach ($row in $sheet2Rows) {
    $ws2.Cells.Item($r, 1).Value = $row[0]
    $ws2.Cells.Item($r, 2).Value = $row[1]
    $ws2.Cells.Item($r, 3).Value = $row[2]
    $ws2.Cells.Item($r, 4).Value = $row[3]
    $ws2.Cells.Item($r, 5).Value = $row[4]
    $r++
}

# Update "6203-6173" sheet (timestamp only - no row/content changes)
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: 06:52:53"
